$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D166").Value = "Variants found: {'Title Association', 'Title Ins', 'Title Insurance'}"
$ws.Range("D167").Value = "Variants found: {'This`nCommitment', 'This Privacy', 'This Notice', 'This Commitment'}"
$ws.Range("D168").Value = "Variants found: {'Commitment`nCondition', 'Commitment Conditions', 'Commitment Condition', 'Commitment Date'}"
$ws.Range("D169").Value = "Variants found: {'Proposed Insured', 'Proposed Amount', 'Proposed Policy'}"
$ws.Range("D170").Value = "Variants found: {'President`nAttest', 'President`nAuthorized'}"
$ws.Range("D172").Value = "Variants found: {'The Policy', 'The West', 'The`nPrivacy', 'The Company', 'The Service', 'The Proposed', 'The Title', 'The State', 'The Land'}"
$ws.Range("D173").Value = "Variants found: {'Real Estate', 'Real Property'}"
$ws.Range("D177").Value = "Variants found: {'Internet Crime', 'Internet Protocol'}"
$ws.Range("D178").Value = "Variants found: {'Privacy`nNotice', 'Privacy Statement', 'Privacy Notice', 'Privacy Inquiry'}"
$ws.Range("D179").Value = "Variants found: {'Personal Information', 'Personal`nInformation'}"
$ws.Range("D181").Value = "Variants found: {'Other Counties', 'Other Sites', 'Other Online'}"
$ws.Range("D184").Value = "Variants found: {'For Virginia', 'For Nevada', 'For Oregon', 'For California'}"

$ws.Range("E165").Value = "LLM enhancement unavailable. Error: Error code: 429 - {'error': {'message': 'You exceeded your current quota, please check your plan and billing details. For more information on this error, read the docs: https://platform.openai.com/docs/guides/error-codes/api-errors.', 'type': 'insufficient_quota', 'param': None, 'code': 'insufficient_quota'}}"
$ws.Range("E198").Value = "LLM error: Error code: 429 - {'error': {'message': 'You exceeded your current quota, please check your plan and billing details. For more information on this error, read the docs: https://platform.openai.com/docs/guides/error-codes/api-errors.', 'type': 'insufficient_quota', 'param': None, 'code': 'insufficient_quota'}}"
